$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + '67.473.20'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'" + '  +2.35%  '
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'" + '3.596.45'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'" + '  +1.51%  '
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'" + '  +0.13%  '
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'" + '197.85'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'" + '  +6.66%  '
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'" + '557.21'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'" + '  -4.20%  '
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'" + '3.593.44'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'" + '  +1.72%  '
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'" + '0.611'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'" + '  +0.35%  '
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'" + '  +0.10%  '
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'" + '0.670'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'" + '  +1.48%  '
$ws.Range("E10").Style = "Normal"
$ws.Range("B11").Value = 'Avalanche'
$ws.Range("C11").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D11").Value = "'" + '56.79'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'" + '  +7.51%  '
$ws.Range("E11").Style = "Normal"
$ws.Range("B12").Value = 'Dogecoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D12").Value = "'" + '0.151'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'" + '  +5.24%  '
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'" + '0.0000289'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'" + '  +14.50%  '
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'" + '9.93'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'" + '  +3.00%  '
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'" + '4.170.44'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'" + '  +1.73%  '
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'" + '3.592.74'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'" + '  +1.56%  '
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'" + '  +0.62%  '
$ws.Range("E17").Style = "Normal"
$ws.Range("B18").Value = 'Chainlink'
$ws.Range("C18").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D18").Value = "'" + '18.79'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'" + '  +3.92%  '
$ws.Range("E18").Style = "Normal"
$ws.Range("B19").Value = 'WrappedBTC'
$ws.Range("C19").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D19").Value = "'" + '67.373.92'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'" + '  +2.63%  '
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'" + '12.24'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'" + '  +1.43%  '
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'" + '1.08'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'" + '  +3.23%  '
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'" + '395.52'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'" + '  +1.15%  '
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'" + '13.13'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'" + '  +26.32%  '
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'" + '  -4.63%  '
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'" + '84.69'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'" + '  +0.40%  '
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'" + '2.93'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'" + '  +2.66%  '
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'" + '12.35'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'" + '  +0.01%  '
$ws.Range("E27").Style = "Normal"
$ws.Range("D29").Value = "'" + '3.81'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'" + '  +8.68%  '
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'" + '8.38'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'" + '  +24.61%  '
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'" + '9.08'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'" + '  +2.97%  '
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'" + '31.35'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'" + '  +2.59%  '
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'" + '676.84'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'" + '  +10.79%  '
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'" + '12.16'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'" + '  +1.01%  '
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'" + '  +3.53%  '
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'" + '63.48'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'" + '  +1.49%  '
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'" + '42.32'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'" + '  +3.56%  '
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'" + '0.433'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'" + '  +17.50%  '
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'" + '  +0.07%  '
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'" + '0.0₃0770'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'" + '  +5.22%  '
$ws.Range("E40").Style = "Normal"
$ws.Range("B41").Value = 'ThetaToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D41").Value = "'" + '3.16'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'" + '  +15.83%  '
$ws.Range("E41").Style = "Normal"
$ws.Range("B42").Value = 'Kaspa'
$ws.Range("C42").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D42").Value = "'" + '0.136'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'" + '  +5.32%  '
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'" + '2.84'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'" + '  +17.18%  '
$ws.Range("E43").Style = "Normal"
$ws.Range("B44").Value = 'Maker'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D44").Value = "'" + '3.199.63'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'" + '  +11.36%  '
$ws.Range("E44").Style = "Normal"
$ws.Range("B45").Value = 'dogwifhat'
$ws.Range("C45").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D45").Value = "'" + '3.03'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'" + '  +34.22%  '
$ws.Range("E45").Style = "Normal"
$ws.Range("B46").Value = 'FirstDigitalUSD'
$ws.Range("C46").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D46").Value = "'" + '0.998'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'" + '  +0.06%  '
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'" + '0.0414'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'" + '  +2.37%  '
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'" + '2.71'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'" + '  +11.08%  '
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'" + '3.11'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'" + '  -0.74%  '
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'" + '  +0.84%  '
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'" + '8.71'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'" + '  +3.55%  '
$ws.Range("E51").Style = "Normal"
